$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom of the sheet upward so row numbers for not-yet-processed
# operations (expressed in terms of the ORIGINAL layout) stay valid.

# 1) Remove the three trailing rows that got relocated/replaced near the bottom
#    (row 334: 005135281/RAFAEL/-0.01, row 335: 004472431/LUIS/-647.19,
#     row 336: 004398253/EULER/-36178.82).
$ws.Rows(336).Delete()
$ws.Rows(335).Delete()
$ws.Rows(334).Delete()

# 2) Insert a new row 005135281/RAFAEL/99.99 right before the existing
#    004463586/MARCIA row (originally row 48).
$ws.Rows(48).Insert()
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "005135281"
$ws.Range("B48").Value = "RAFAEL"
$ws.Range("C48").Value = 99.99

# 3) Insert a new row 004853111/MARCONDES/3278.26 right before the existing
#    005654122/ELIANE row (originally row 15).
$ws.Rows(15).Insert()
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "004853111"
$ws.Range("B15").Value = "MARCONDES"
$ws.Range("C15").Value = 3278.26

# 4) Remove the 008032257/SARA/19006 row (row 9).
$ws.Rows(9).Delete()
